$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.893.08'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.884.56'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.53'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.86%  '

$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4621'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.69%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4090'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.78%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.28'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07985'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.67%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9876'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.66'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.918.55'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.902'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.060'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.03%  '

$ws.Range("E16").Value = '  +0.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.89'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001028'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06567'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.44'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.879.68'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.397'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.21'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.204'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.086.26'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.36'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.61'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.098'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.385'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.77'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9734'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09341'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.95%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.603'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.41%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.400'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.269'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06047'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02228'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.251'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("E40").Value = '  +0.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5761'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.87%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1819'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.91%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.10'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.251'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.262'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +9.30%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5455'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.84%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.90'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.898'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.61%  '

$ws.Range("E50").Value = '  -6.13%  '

$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.67'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +14.73%  '
